$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(
    5980,5950,5930,5900,5880,5850,5840,5830,5820,5830,5830,5850,5890,5910,5940,5980,6030,6100,6190,6290,6380,6510,6640,6770,6980,7080,7150,7180,7150,7110,7060,6990,6890,6790,6670,6530,6390,6260,6140,6020,5890,5790,5700,5630,5550,5500,5460,5440,5440,5440,5440,5460,5470,5490,5510,5540,5600,5650,5710,5780,5880,6000,6140,6300,6460,6640,6820,7000,7240,7420,7590,7750,7900,8000,8050,8040,7990,7940,7870,7800,7680,7560,7450,7310,7120,6970,6800,6650,6530,6400,6280,6170,6130,6080,6020,5960
)

$bValues = @(
    45736,45736.01041666666,45736.02083333334,45736.03125,45736.04166666666,45736.05208333334,45736.0625,45736.07291666666,45736.08333333334,45736.09375,45736.10416666666,45736.11458333334,45736.125,45736.13541666666,45736.14583333334,45736.15625,45736.16666666666,45736.17708333334,45736.1875,45736.19791666666,45736.20833333334,45736.21875,45736.22916666666,45736.23958333334,45736.25,45736.26041666666,45736.27083333334,45736.28125,45736.29166666666,45736.30208333334,45736.3125,45736.32291666666,45736.33333333334,45736.34375,45736.35416666666,45736.36458333334,45736.375,45736.38541666666,45736.39583333334,45736.40625,45736.41666666666,45736.42708333334,45736.4375,45736.44791666666,45736.45833333334,45736.46875,45736.47916666666,45736.48958333334,45736.5,45736.51041666666,45736.52083333334,45736.53125,45736.54166666666,45736.55208333334,45736.5625,45736.57291666666,45736.58333333334,45736.59375,45736.60416666666,45736.61458333334,45736.625,45736.63541666666,45736.64583333334,45736.65625,45736.66666666666,45736.67708333334,45736.6875,45736.69791666666,45736.70833333334,45736.71875,45736.72916666666,45736.73958333334,45736.75,45736.76041666666,45736.77083333334,45736.78125,45736.79166666666,45736.80208333334,45736.8125,45736.82291666666,45736.83333333334,45736.84375,45736.85416666666,45736.86458333334,45736.875,45736.88541666666,45736.89583333334,45736.90625,45736.91666666666,45736.92708333334,45736.9375,45736.94791666666,45736.95833333334,45736.96875,45736.97916666666,45736.98958333334
)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
